# Refined metadata to be additional tab
#
# 1) Update the two panel-query timestamps on the "data" sheet.
# 2) Add a new "metadata" sheet (placed after "data") describing the
#    PanelApp query that produced the data sheet.

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item(1)

# --- 1. refresh timestamps on the "data" sheet --------------------------
$data.Range("F2").Value = "2021-10-05 14:33:51.446028"
$data.Range("F3").Value = "2021-10-05 14:33:51.446036"

# --- 2. add the "metadata" sheet -----------------------------------------
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# header row (row 1) ------------------------------------------------------
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# copy the header style (bold / bordered / centered) from the "data" sheet
$data.Range("B1:F1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# data row (row 2) ---------------------------------------------------------
$meta.Range("A2").Value = 0
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$meta.Range("B2").Value = "Foveal Hypoplasia"
$meta.Range("C2").Value = 3150
$meta.Range("D2").Value = "'0.7"
$meta.Range("E2").Value = "2020-12-03T04:43:40.968133Z"
$meta.Range("F2").Value = "2021-10-05 14:33:51.442119"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3150/?format=json"

# D2 must stay plain text ("0.7") without picking up a quote-prefix style;
# re-paste formatting from an untouched, default-styled cell to clear it.
$data.Range("H10").Copy()
$meta.Range("D2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$data.Activate()
$data.Range("A1").Select()
